$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 19189.875
$ws.Range("I21").Value = 23367.092
$ws.Range("K21").Value = 23367.092
$ws.Range("M21").Value = -22899.092

$ws.Range("H23").Value = 19189.875
$ws.Range("I23").Value = 23367.092
$ws.Range("K23").Value = 23367.092
$ws.Range("M23").Value = -23133.092

$ws.Range("H28").Value = 606.44446
$ws.Range("I28").Value = 661.4167
$ws.Range("J28").Value = 166.66667
$ws.Range("K28").Value = 661.4167
$ws.Range("L28").Value = 166.66667
$ws.Range("M28").Value = -176.4167
$ws.Range("N28").Value = -1136.66667

$ws.Range("H62").Value = 4198.3335
$ws.Range("I62").Value = 4886.6665
$ws.Range("J62").Value = 2133.3333
$ws.Range("K62").Value = 4886.6665
$ws.Range("L62").Value = 2133.3333
$ws.Range("M62").Value = -4262.6665
$ws.Range("N62").Value = -3381.3333

$ws.Range("H65").Value = 4198.3335
$ws.Range("I65").Value = 4886.6665
$ws.Range("J65").Value = 2133.3333
$ws.Range("K65").Value = 24433.3325
$ws.Range("L65").Value = 10666.6665
$ws.Range("M65").Value = -21313.3325
$ws.Range("N65").Value = -16906.6665

$ws.Range("H92").Value = 446
$ws.Range("I92").Value = 389.45834
$ws.Range("J92").Value = 1124.5
$ws.Range("K92").Value = 389.45834
$ws.Range("L92").Value = 1124.5
$ws.Range("M92").Value = 858.54166
$ws.Range("N92").Value = -3620.5

$ws.Range("H95").Value = 33721.5
$ws.Range("J95").Value = 33721.5
$ws.Range("L95").Value = 33721.5
$ws.Range("N95").Value = -39213.5

$ws.Range("H113").Value = 2500
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 2500
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 2500
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -9008

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3320.4
$ws.Range("I32").Value = 3320.4
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 3320.4
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -3033.4
$ws.Range("N32").ClearContents()

$ws.Range("H97").Value = 946.9
$ws.Range("I97").Value = 642.7646999999999
$ws.Range("J97").Value = 2670.3333
$ws.Range("K97").Value = 642.7646999999999
$ws.Range("L97").Value = 2670.3333
$ws.Range("M97").Value = -146.7646999999999
$ws.Range("N97").Value = -3662.3333

$ws.Range("H122").Value = 971
$ws.Range("I122").Value = 971
$ws.Range("K122").Value = 2913
$ws.Range("M122").Value = -463

$ws.Range("H132").Value = 8476884
$ws.Range("I132").Value = 16668234
$ws.Range("J132").Value = 3073.9312
$ws.Range("K132").Value = 50004702
$ws.Range("L132").Value = 9221.793600000001
$ws.Range("M132").Value = -50002172
$ws.Range("N132").Value = -14281.7936

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 960
$ws.Range("I94").Value = 820
$ws.Range("J94").Value = 1100
$ws.Range("K94").Value = 820
$ws.Range("L94").Value = 1100
$ws.Range("M94").Value = -369
$ws.Range("N94").Value = -2002

$ws.Range("H97").Value = 6960.857
$ws.Range("I97").Value = 1745.2
$ws.Range("J97").Value = 20000
$ws.Range("K97").Value = 1745.2
$ws.Range("L97").Value = 20000
$ws.Range("M97").Value = -754.2
$ws.Range("N97").Value = -21982

$ws.Range("H99").Value = 1561.3889
$ws.Range("I99").Value = 1504.4073
$ws.Range("K99").Value = 1504.4073
$ws.Range("M99").Value = -6.407300000000077

$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

$ws.Range("H105").Value = 2103.6943
$ws.Range("I105").Value = 1963.1904
$ws.Range("K105").Value = 1963.1904
$ws.Range("M105").Value = -216.1904

$ws.Range("H107").Value = 1941.3334
$ws.Range("I107").Value = 1882.4
$ws.Range("J107").Value = 2059.2
$ws.Range("K107").Value = 1882.4
$ws.Range("L107").Value = 2059.2
$ws.Range("M107").Value = 37.59999999999991
$ws.Range("N107").Value = -5899.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8956.92
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 8956.92
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 8956.92
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -9546.92

$ws.Range("H33").Value = 11982.75
$ws.Range("I33").Value = 11982.75
$ws.Range("K33").Value = 11982.75
$ws.Range("M33").Value = -11603.75

$ws.Range("H34").Value = 8956.92
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 8956.92
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 8956.92
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -9360.92

$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()

$ws.Range("H99").Value = 2008.2307
$ws.Range("I99").Value = 2100.7778
$ws.Range("J99").Value = 1800
$ws.Range("K99").Value = 2100.7778
$ws.Range("L99").Value = 1800
$ws.Range("M99").Value = -602.7777999999998
$ws.Range("N99").Value = -4796

$ws.Range("H126").Value = 2008.2307
$ws.Range("I126").Value = 2100.7778
$ws.Range("J126").Value = 1800
$ws.Range("K126").Value = 6302.3334
$ws.Range("L126").Value = 5400
$ws.Range("M126").Value = -3832.3334
$ws.Range("N126").Value = -10340

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 205.72728
$ws.Range("I12").Value = 217.18182
$ws.Range("J12").Value = 194.27272
$ws.Range("K12").Value = 651.5454599999999
$ws.Range("L12").Value = 582.81816
$ws.Range("M12").Value = -478.5454599999999
$ws.Range("N12").Value = -928.81816

$ws.Range("H68").Value = 595.3333
$ws.Range("I68").Value = 300
$ws.Range("J68").Value = 743
$ws.Range("K68").Value = 900
$ws.Range("L68").Value = 2229
$ws.Range("M68").Value = -89
$ws.Range("N68").Value = -3851

$ws.Range("H71").Value = 595.3333
$ws.Range("I71").Value = 300
$ws.Range("J71").Value = 743
$ws.Range("K71").Value = 2700
$ws.Range("L71").Value = 6687
$ws.Range("M71").Value = 1356
$ws.Range("N71").Value = -14799

$ws.Range("H113").Value = 3751.0938
$ws.Range("I113").Value = 6877.125
$ws.Range("J113").Value = 625.0625
$ws.Range("K113").Value = 20631.375
$ws.Range("L113").Value = 1875.1875
$ws.Range("M113").Value = -18461.375
$ws.Range("N113").Value = -6215.1875

$ws.Range("H137").Value = 40007840
$ws.Range("I137").Value = 3070.6667
$ws.Range("J137").Value = 76935320
$ws.Range("K137").Value = 9212.000100000001
$ws.Range("L137").Value = 230805960
$ws.Range("M137").Value = -4112.000100000001
$ws.Range("N137").Value = -230816160

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H101").Value = 48996
$ws.Range("J101").Value = 48996
$ws.Range("L101").Value = 48996
$ws.Range("N101").Value = -55486

$ws.Range("H132").Value = 2453.1282
$ws.Range("I132").Value = 1798.75
$ws.Range("J132").Value = 4118.8184
$ws.Range("K132").Value = 5396.25
$ws.Range("L132").Value = 12356.4552
$ws.Range("M132").Value = -2866.25
$ws.Range("N132").Value = -17416.4552

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1676.2222
$ws.Range("I100").Value = 1405.1428
$ws.Range("K100").Value = 1405.1428
$ws.Range("M100").Value = -864.1428000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 41390.727
$ws.Range("J135").Value = 41390.727
$ws.Range("L135").Value = 41390.727
$ws.Range("N135").Value = -51530.727
